$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new (hidden-less) autofilter-related defined name, mirroring the
# existing chain of _xlnm._FilterDatabase* names already on the sheet.
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0", "=Sheet1!`$A`$1:`$J`$41")

# --- key_json value shared by the four new rows
$keyJson = '{"store_att_name_1": "store_type", "store_att_value_1": "Supers A"}'

$dataJsonGumCheckout = '{"Target": "", "Weight": 20, "kpi_child": "", "KPI Family": "Linear SOS", "kpi_parent": "", "score_logic": "Tiered", "Template Name": "Checkout Gum & Confectionary", "KPI Level 2 Name": "Gum & Fruity", "score_cond_score_1": 0, "score_cond_score_2": 0.75, "score_cond_score_3": 1, "score_cond_score_4": "", "score_cond_target_1": 0.6, "score_cond_target_2": 0.7, "score_cond_target_3": 1.01, "score_cond_target_4": "", "exclude_param_type_1": "sub_category_fk", "exclude_param_type_2": "product_type", "exclude_param_value_1": 20, "exclude_param_value_2": ["Other", "Empty"], "param_type_2/denom_type": "category_fk", "param_value_2/denom_value": 10, "param_type_1/numerator_type": "manufacturer_fk", "param_value_1/numerator_value": 762, "exclude_param_2_exception_type": "brand_fk", "exclude_param_2_exception_value": 82}'

$dataJsonGumMain = '{"Target": "", "Weight": 15, "kpi_child": "", "KPI Family": "Linear SOS", "kpi_parent": "", "score_logic": "Tiered", "Template Name": "Main Shelf Gum and Confectionary", "KPI Level 2 Name": "Gum & Fruity", "score_cond_score_1": 0, "score_cond_score_2": 0.8, "score_cond_score_3": 0.9, "score_cond_score_4": 1, "score_cond_target_1": 0.58, "score_cond_target_2": 0.6, "score_cond_target_3": 0.62, "score_cond_target_4": 1.01, "exclude_param_type_1": "sub_category_fk", "exclude_param_type_2": "product_type", "exclude_param_value_1": 20, "exclude_param_value_2": ["Other", "Empty"], "param_type_2/denom_type": "category_fk", "param_value_2/denom_value": 10, "param_type_1/numerator_type": "manufacturer_fk", "param_value_1/numerator_value": 762, "exclude_param_2_exception_type": "brand_fk", "exclude_param_2_exception_value": 82}'

# These two re-use the (pre-existing) chocolate tiered/binary data_json blobs
$dataJsonChocolateMain = '{"Target": "", "Weight": 25, "kpi_child": "", "KPI Family": "Linear SOS", "kpi_parent": "", "score_logic": "Tiered", "Template Name": "Main shelf chocolate", "KPI Level 2 Name": "Chocolate & Ice Cream", "score_cond_score_1": 1, "score_cond_score_2": 0.9, "score_cond_score_3": 0.8, "score_cond_score_4": "", "score_cond_target_1": 0.4, "score_cond_target_2": 0.37, "score_cond_target_3": 0.35, "score_cond_target_4": "", "exclude_param_type_1": "", "exclude_param_value_1": "", "param_type_2/denom_type": "category_fk", "param_value_2/denom_value": 6, "param_type_1/numerator_type": "manufacturer_fk", "param_value_1/numerator_value": 2}'

$dataJsonChocolateCheckout = '{"Target": 0.7, "Weight": 5, "kpi_child": "", "KPI Family": "Linear SOS", "kpi_parent": "", "score_logic": "Binary", "Template Name": "Checkout Chocolate", "KPI Level 2 Name": "Chocolate & Ice Cream", "score_cond_score_1": "", "score_cond_score_2": "", "score_cond_score_3": "", "score_cond_score_4": "", "score_cond_target_1": "", "score_cond_target_2": "", "score_cond_target_3": "", "score_cond_target_4": "", "exclude_param_type_1": "", "exclude_param_value_1": "", "param_type_2/denom_type": "category_fk", "param_value_2/denom_value": 6, "param_type_1/numerator_type": "manufacturer_fk", "param_value_1/numerator_value": 2}'

# --- New rows 42-45
# Row 42
$ws.Range("A42").Value = 83
$ws.Range("B42").Value = 2
$ws.Range("C42").Value = 3031
$ws.Range("D42").Value = $keyJson
$ws.Range("E42").Value = $dataJsonGumCheckout
$ws.Range("F42").Value = 43468
$ws.Range("F42").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("H42").Value = 43677.5600115741
$ws.Range("H42").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("I42").Value = "atomic_level"
$ws.Range("J42").Value = "SOS - Gum Checkout"

# Row 43
$ws.Range("A43").Value = 85
$ws.Range("B43").Value = 2
$ws.Range("C43").Value = 3032
$ws.Range("D43").Value = $keyJson
$ws.Range("E43").Value = $dataJsonGumMain
$ws.Range("F43").Value = 43468
$ws.Range("F43").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("H43").Value = 43677.5600115741
$ws.Range("H43").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("I43").Value = "atomic_level"
$ws.Range("J43").Value = "SOS - Gum Main"

# Row 44
$ws.Range("A44").Value = 86
$ws.Range("B44").Value = 2
$ws.Range("C44").Value = 3030
$ws.Range("D44").Value = $keyJson
$ws.Range("E44").Value = $dataJsonChocolateMain
$ws.Range("F44").Value = 43468
$ws.Range("F44").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("H44").Value = 43677.5600115741
$ws.Range("H44").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("I44").Value = "atomic_level"
$ws.Range("J44").Value = "SOS - Chocolate Main"

# Row 45
$ws.Range("A45").Value = 87
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = 3029
$ws.Range("D45").Value = $keyJson
$ws.Range("E45").Value = $dataJsonChocolateCheckout
$ws.Range("F45").Value = 43468
$ws.Range("F45").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("H45").Value = 43677.5600115741
$ws.Range("H45").NumberFormat = "YYYY\-MM\-DD\ HH:MM:SS"
$ws.Range("I45").Value = "atomic_level"
$ws.Range("J45").Value = "SOS - Chocolate Checkout"

# --- Column width tweaks to match the resulting (slightly narrower/wider)
# layout; values are chosen so the pixel-quantized COM ColumnWidth lands as
# close as possible to the author's original fractional widths.
$ws.Columns.Item(2).ColumnWidth = 9.5
$ws.Columns.Item(3).ColumnWidth = 12.166667
$ws.Columns.Item(4).ColumnWidth = 58.833333
$ws.Columns.Item(5).ColumnWidth = 144.0
$ws.Columns.Item(6).ColumnWidth = 20.333333
$ws.Columns.Item(7).ColumnWidth = 9.5
$ws.Columns.Item(8).ColumnWidth = 19.166667
$ws.Columns.Item(9).ColumnWidth = 34.0

# --- Selection follows the last edited row/column, as in the source edit
$ws.Range("E44").Select()
